$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 192, shifting existing rows 192:226 down to 193:227
$ws.Rows.Item(192).Insert()

# Populate the newly inserted row 192 with the new price-record data
$ws.Range("A192").Value = 5
$ws.Range("B192").Value = "Macroferia Regional de Talca"
$ws.Range("C192").Value = "Maule"
$ws.Range("D192").Value = 44522
$ws.Range("E192").Value = 7
$ws.Range("F192").Value = 100114013
$ws.Range("G192").Value = "Zanahoria"
$ws.Range("H192").Value = "Sin especificar"
$ws.Range("I192").Value = "Primera"
$ws.Range("J192").Value = 400
$ws.Range("K192").Value = 9000
$ws.Range("L192").Value = 9000
$ws.Range("M192").Value = 9000
$ws.Range("N192").Value = "$/saco 20 kilos"
$ws.Range("O192").Value = "Provincia del Elquí"
$ws.Range("P192").Value = 450
$ws.Range("Q192").Value = 20
$ws.Range("R192").Value = "Hortaliza"
